$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 829.3674
$ws.Range("J17").Value = 829.3674
$ws.Range("L17").Value = 2488.1022
$ws.Range("N17").Value = -2824.1022

$ws.Range("H112").Value = 2865.25
$ws.Range("J112").Value = 3054.8276
$ws.Range("L112").Value = 9164.4828
$ws.Range("N112").Value = -11380.4828

$ws.Range("H127").Value = 1500.8889
$ws.Range("I127").Value = 343.6
$ws.Range("J127").Value = 2947.5
$ws.Range("K127").Value = 1030.8
$ws.Range("L127").Value = 8842.5
$ws.Range("M127").Value = 3929.2
$ws.Range("N127").Value = -18762.5

$ws.Range("H132").Value = 10759317
$ws.Range("I132").Value = 12352920
$ws.Range("J132").Value = 2490
$ws.Range("K132").Value = 37058760
$ws.Range("L132").Value = 7470
$ws.Range("M132").Value = -37056230
$ws.Range("N132").Value = -12530

$ws.Range("H138").Value = 1400.53
$ws.Range("I138").Value = 900.1212
$ws.Range("J138").Value = 1647
$ws.Range("K138").Value = 2700.3636
$ws.Range("L138").Value = 4941
$ws.Range("M138").Value = 2439.6364
$ws.Range("N138").Value = -15221

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1237.091
$ws.Range("I45").Value = 1350.25
$ws.Range("K45").Value = 1350.25
$ws.Range("M45").Value = -973.25

$ws.Range("H110").Value = 1086.8077
$ws.Range("I110").Value = 634.6
$ws.Range("J110").Value = 1703.4546
$ws.Range("K110").Value = 634.6
$ws.Range("L110").Value = 1703.4546
$ws.Range("M110").Value = 1410.4
$ws.Range("N110").Value = -5793.4546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 252474220
$ws.Range("I105").Value = 336631300
$ws.Range("K105").Value = 336631300
$ws.Range("M105").Value = -336629553

$ws.Range("H107").Value = 1712.1538
$ws.Range("I107").Value = 1241.5
$ws.Range("K107").Value = 1241.5
$ws.Range("M107").Value = 678.5

$ws.Range("H134").Value = 5343.3335
$ws.Range("I134").Value = 845.875
$ws.Range("J134").Value = 23333.166
$ws.Range("K134").Value = 2537.625
$ws.Range("L134").Value = 69999.49800000001
$ws.Range("M134").Value = -2.625
$ws.Range("N134").Value = -75069.49800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1560.8928
$ws.Range("I31").Value = 1373.409
$ws.Range("J31").Value = 2248.3333
$ws.Range("K31").Value = 1373.409
$ws.Range("L31").Value = 2248.3333
$ws.Range("M31").Value = -1078.409
$ws.Range("N31").Value = -2838.3333

$ws.Range("H34").Value = 1560.8928
$ws.Range("I34").Value = 1373.409
$ws.Range("J34").Value = 2248.3333
$ws.Range("K34").Value = 1373.409
$ws.Range("L34").Value = 2248.3333
$ws.Range("M34").Value = -1171.409
$ws.Range("N34").Value = -2652.3333

$ws.Range("H62").Value = 66668668
$ws.Range("I62").Value = 3000
$ws.Range("K62").Value = 3000
$ws.Range("M62").Value = -2376

$ws.Range("H65").Value = 66668668
$ws.Range("I65").Value = 3000
$ws.Range("K65").Value = 15000
$ws.Range("M65").Value = -11880

$ws.Range("H122").Value = 950.5
$ws.Range("I122").Value = 797.625
$ws.Range("K122").Value = 2392.875
$ws.Range("M122").Value = 57.125

$ws.Range("H132").Value = 4371.375
$ws.Range("I132").Value = 5249.24
$ws.Range("J132").Value = 2908.2666
$ws.Range("K132").Value = 15747.72
$ws.Range("L132").Value = 8724.799800000001
$ws.Range("M132").Value = -13217.72
$ws.Range("N132").Value = -13784.7998

$ws.Range("H134").Value = 1513.0975
$ws.Range("I134").Value = 1501.0938
$ws.Range("J134").Value = 1555.7778
$ws.Range("K134").Value = 4503.2814
$ws.Range("L134").Value = 4667.3334
$ws.Range("M134").Value = -1968.2814
$ws.Range("N134").Value = -9737.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 68.76000000000001
$ws.Range("I12").Value = 77.125
$ws.Range("J12").Value = 64.82353000000001
$ws.Range("K12").Value = 231.375
$ws.Range("L12").Value = 194.47059
$ws.Range("M12").Value = -58.375
$ws.Range("N12").Value = -540.47059

$ws.Range("H18").Value = 1222.2222
$ws.Range("I18").Value = 470
$ws.Range("K18").Value = 1410
$ws.Range("M18").Value = -1241

$ws.Range("H69").Value = 2000.56
$ws.Range("I69").Value = 1800
$ws.Range("J69").Value = 2018
$ws.Range("K69").Value = 5400
$ws.Range("L69").Value = 6054
$ws.Range("M69").Value = -4589
$ws.Range("N69").Value = -7676

$ws.Range("H72").Value = 2000.56
$ws.Range("I72").Value = 1800
$ws.Range("J72").Value = 2018
$ws.Range("K72").Value = 16200
$ws.Range("L72").Value = 18162
$ws.Range("M72").Value = -12144
$ws.Range("N72").Value = -26274

$ws.Range("H92").Value = 662.0909
$ws.Range("I92").Value = 685.375
$ws.Range("K92").Value = 2056.125
$ws.Range("M92").Value = -808.125

$ws.Range("H103").Value = 1778.3334
$ws.Range("I103").Value = 561
$ws.Range("J103").Value = 3300
$ws.Range("K103").Value = 1683
$ws.Range("L103").Value = 9900
$ws.Range("M103").Value = -804
$ws.Range("N103").Value = -11658

$ws.Range("H104").Value = 4566.8335
$ws.Range("J104").Value = 5433.3335
$ws.Range("L104").Value = 16300.0005
$ws.Range("N104").Value = -21542.0005

$ws.Range("H121").Value = 707.61536
$ws.Range("J121").Value = 811.6667
$ws.Range("L121").Value = 2435.0001
$ws.Range("N121").Value = -5055.0001

$ws.Range("H131").Value = 15152709
$ws.Range("J131").Value = 1266.3442
$ws.Range("L131").Value = 3799.0326
$ws.Range("N131").Value = -13879.0326

$ws.Range("H132").Value = 782.5
$ws.Range("I132").Value = 270
$ws.Range("K132").Value = 2430
$ws.Range("M132").Value = 100

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2197.2144
$ws.Range("I102").Value = 2213.4167
$ws.Range("K102").Value = 2213.4167
$ws.Range("M102").Value = -591.4167000000002

$ws.Range("H126").Value = 1998.2632
$ws.Range("I126").Value = 1861.9286
$ws.Range("K126").Value = 5585.7858
$ws.Range("M126").Value = -3115.7858

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 201.65
$ws.Range("I55").Value = 110.92857
$ws.Range("J55").Value = 413.33334
$ws.Range("K55").Value = 110.92857
$ws.Range("L55").Value = 413.33334
$ws.Range("M55").Value = 62.07143000000001
$ws.Range("N55").Value = -759.33334

$ws.Range("H132").Value = 23002.469
$ws.Range("I132").Value = 1426.6957
$ws.Range("J132").Value = 43679.25
$ws.Range("K132").Value = 4280.0871
$ws.Range("L132").Value = 131037.75
$ws.Range("M132").Value = -1750.0871
$ws.Range("N132").Value = -136097.75

$ws.Range("H136").Value = 10008.333
$ws.Range("I136").Value = 10645.454
$ws.Range("K136").Value = 31936.362
$ws.Range("M136").Value = -29386.362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 17900
$ws.Range("J64").Value = 17900
$ws.Range("L64").Value = 17900
$ws.Range("N64").Value = -18396

$ws.Range("H67").Value = 17900
$ws.Range("J67").Value = 17900
$ws.Range("L67").Value = 17900
$ws.Range("N67").Value = -19616

$ws.Range("H122").Value = 9631891
$ws.Range("I122").Value = 11306673
$ws.Range("K122").Value = 33920019
$ws.Range("M122").Value = -33917569

$ws.Range("H126").Value = 52632640
$ws.Range("I126").Value = 71429560
$ws.Range("J126").Value = 1262.4
$ws.Range("K126").Value = 214288680
$ws.Range("L126").Value = 3787.2
$ws.Range("M126").Value = -214286210
$ws.Range("N126").Value = -8727.200000000001

$ws.Range("H136").Value = 825.2857
$ws.Range("I136").Value = 825.2857
$ws.Range("K136").Value = 2475.8571
$ws.Range("M136").Value = 74.14289999999983
